$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (keeps
#    all formatting/styles identical) and inserting it right before
#    "2022-Q2" in the tab order.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$origIndex = $q2.Index
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item($origIndex)
$q3.Name = "2022-Q3"

# Make sure the text-like numeric columns (B..G) keep their original
# "stored as text" semantics (leading zeros / formatted numbers) when
# we overwrite them below.
$q3.Range("B2:G4").NumberFormat = "@"

# Row 2 - first fund entry: code (160910) and H2 (rank) are unchanged,
# only name / size / position figures move.
$q3.Range("C2").Value = "大成创新成长混合（LOF）"
$q3.Range("D2").Value = "10.78"
$q3.Range("E2").Value = "85.81"
$q3.Range("F2").Value = "7.39"
$q3.Range("G2").Value = "0.7966"

# Row 3 - second fund entry is replaced entirely.
$q3.Range("B3").Value = "015707"
$q3.Range("C3").Value = "安信新能源主题股票A"
$q3.Range("D3").Value = "0.18"
$q3.Range("E3").Value = "53.26"
$q3.Range("F3").Value = "1.87"
$q3.Range("G3").Value = "0.0034"
$q3.Range("H3").Value = 10

# Row 4 - third fund entry is replaced entirely.
$q3.Range("B4").Value = "015708"
$q3.Range("C4").Value = "安信新能源主题股票C"
$q3.Range("D4").Value = "0.14"
$q3.Range("E4").Value = "53.26"
$q3.Range("F4").Value = "1.87"
$q3.Range("G4").Value = "0.0026"
$q3.Range("H4").Value = 10

# ------------------------------------------------------------------
# 2. Update the "总计" (Total) summary sheet: add a row for the new
#    2022-Q3 quarter above the existing history, shifting the older
#    quarters down by one row.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the old "2022-Q1" row (row 3) down to row 4.
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
# Duplicate the old "2022-Q2" row (row 2) into row 3 - it keeps its values.
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# Fix up the running index column and row labels.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Turn row 2 into the new 2022-Q3 entry.
$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.8

# Restore the original active sheet/selection (总计 was the workbook's
# active tab before this edit).
$total.Activate()
$null = $total.Range("A1").Select()
